$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" '68.465.02'
Set-TextCell $ws "E2" '  +0.19%  '
Set-TextCell $ws "D3" '3.853.08'
Set-TextCell $ws "E3" '  -1.49%  '
Set-TextCell $ws "E4" '  +0.07%  '
Set-TextCell $ws "D5" '522.39'
Set-TextCell $ws "E5" '  +7.17%  '
Set-TextCell $ws "D6" '142.20'
Set-TextCell $ws "E6" '  -2.46%  '
Set-TextCell $ws "D7" '0.602'
Set-TextCell $ws "E7" '  -3.02%  '
Set-TextCell $ws "D8" '1.00'
Set-TextCell $ws "E8" '  +0.17%  '
Set-TextCell $ws "D9" '0.709'
Set-TextCell $ws "E9" '  -4.27%  '
Set-TextCell $ws "D10" '0.166'
Set-TextCell $ws "E10" '  -6.58%  '
Set-TextCell $ws "D11" '0.0000319'
Set-TextCell $ws "E11" '  -7.94%  '
Set-TextCell $ws "D12" '41.43'
Set-TextCell $ws "E12" '  -3.51%  '
Set-TextCell $ws "B13" 'WrappedliquidstakedEther2.0'
Set-TextCell $ws "C13" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell $ws "D13" '4.475.14'
Set-TextCell $ws "E13" '  -1.32%  '
Set-TextCell $ws "B14" 'Polkadot'
Set-TextCell $ws "C14" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws "D14" '10.15'
Set-TextCell $ws "E14" '  -3.35%  '
Set-TextCell $ws "D15" '21.81'
Set-TextCell $ws "E15" '  +9.02%  '
Set-TextCell $ws "D16" '3.847.69'
Set-TextCell $ws "E16" '  -1.97%  '
Set-TextCell $ws "D17" '14.03'
Set-TextCell $ws "E17" '  -0.83%  '
Set-TextCell $ws "B18" 'TRON'
Set-TextCell $ws "C18" 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws "D18" '0.134'
Set-TextCell $ws "E18" '  -1.31%  '
Set-TextCell $ws "B19" 'Polygon'
Set-TextCell $ws "C19" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell $ws "D19" '1.21'
Set-TextCell $ws "E19" '  +4.31%  '
Set-TextCell $ws "D20" '68.562.32'
Set-TextCell $ws "E20" '  +0.18%  '
Set-TextCell $ws "D21" '415.17'
Set-TextCell $ws "E21" '  -4.06%  '
Set-TextCell $ws "D22" '3.38'
Set-TextCell $ws "E22" '  -6.10%  '
Set-TextCell $ws "D23" '13.93'
Set-TextCell $ws "E23" '  -5.90%  '
Set-TextCell $ws "D24" '86.59'
Set-TextCell $ws "E24" '  -3.71%  '
Set-TextCell $ws "D25" '3.92'
Set-TextCell $ws "E25" '  +4.93%  '
Set-TextCell $ws "D26" '11.15'
Set-TextCell $ws "E26" '  -8.88%  '
Set-TextCell $ws "D27" '10.47'
Set-TextCell $ws "E27" '  -6.08%  '
Set-TextCell $ws "D28" '35.15'
Set-TextCell $ws "E28" '  -5.54%  '
Set-TextCell $ws "D29" '682.42'
Set-TextCell $ws "E29" '  -4.45%  '
Set-TextCell $ws "D30" '13.02'
Set-TextCell $ws "E30" '  -2.67%  '
Set-TextCell $ws "D31" '0.124'
Set-TextCell $ws "E31" '  -5.93%  '
Set-TextCell $ws "D32" '2.80'
Set-TextCell $ws "E32" '  -3.81%  '
Set-TextCell $ws "D33" '65.60'
Set-TextCell $ws "E33" '  +6.96%  '
Set-TextCell $ws "B34" 'InjectiveProtocol'
Set-TextCell $ws "C34" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws "D34" '39.57'
Set-TextCell $ws "E34" '  -3.66%  '
Set-TextCell $ws "D35" '5.81'
Set-TextCell $ws "E35" '  -5.09%  '
Set-TextCell $ws "B36" 'TheGraph'
Set-TextCell $ws "C36" 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell $ws "D36" '0.419'
Set-TextCell $ws "E36" '  -11.49%  '
Set-TextCell $ws "B37" 'PEPE'
Set-TextCell $ws "C37" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws "D37" '0.0₃0829'
Set-TextCell $ws "E37" '  -4.58%  '
Set-TextCell $ws "D38" '0.996'
Set-TextCell $ws "E38" '  -0.20%  '
Set-TextCell $ws "D39" '0.147'
Set-TextCell $ws "E39" '  -1.21%  '
Set-TextCell $ws "D40" '1.00'
Set-TextCell $ws "E40" '  +0.10%  '
Set-TextCell $ws "D41" '3.20'
Set-TextCell $ws "E41" '  +2.95%  '
Set-TextCell $ws "D42" '0.0472'
Set-TextCell $ws "E42" '  -4.71%  '
Set-TextCell $ws "D43" '3.14'
Set-TextCell $ws "E43" '  +5.47%  '
Set-TextCell $ws "D44" '2.74'
Set-TextCell $ws "E44" '  -9.98%  '
Set-TextCell $ws "D45" '3.38'
Set-TextCell $ws "E45" '  -1.25%  '
Set-TextCell $ws "D46" '0.138'
Set-TextCell $ws "E46" '  -2.93%  '
Set-TextCell $ws "B47" 'Stacks'
Set-TextCell $ws "C47" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws "D47" '2.90'
Set-TextCell $ws "E47" '  +2.92%  '
Set-TextCell $ws "B48" 'Monero'
Set-TextCell $ws "C48" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws "D48" '143.44'
Set-TextCell $ws "E48" '  +0.91%  '
Set-TextCell $ws "B49" 'Maker'
Set-TextCell $ws "C49" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws "D49" '2.694.75'
Set-TextCell $ws "E49" '  +10.73%  '
Set-TextCell $ws "D50" '0.0₆0336'
Set-TextCell $ws "E50" '  -9.60%  '
Set-TextCell $ws "B51" 'LidoDAOToken'
Set-TextCell $ws "C51" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws "D51" '3.24'
Set-TextCell $ws "E51" '  -4.32%  '
